$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price/Volume columns to Text so numeric-looking strings
# (e.g. "605.17", "1.00") are not auto-converted to numbers by Excel.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

# Row 2 - Bitcoin
$ws.Range("D2").Value = "64.325.74"
$ws.Range("E2").Value = "  +1.99%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.643.63"
$ws.Range("E3").Value = "  +0.75%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.00%  "

# Row 5 - BNB
$ws.Range("D5").Value = "605.17"
$ws.Range("E5").Value = "  +0.16%  "

# Row 6 - Solana
$ws.Range("D6").Value = "151.99"
$ws.Range("E6").Value = "  +3.91%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.01%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  +1.25%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  +1.79%  "

# Row 10 - Cardano
$ws.Range("D10").Value = "0.391"
$ws.Range("E10").Value = "  +8.05%  "

# Row 11 - Toncoin
$ws.Range("D11").Value = "5.68"
$ws.Range("E11").Value = "  +1.24%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  -0.65%  "

# Row 13 - Avalanche
$ws.Range("D13").Value = "27.79"
$ws.Range("E13").Value = "  +2.25%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "3.119.40"
$ws.Range("E14").Value = "  +0.90%  "

# Row 15 - WrappedBTC
$ws.Range("D15").Value = "64.151.75"
$ws.Range("E15").Value = "  +1.95%  "

# Row 16 - ShibaInu
$ws.Range("E16").Value = "  +3.22%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "2.645.88"
$ws.Range("E17").Value = "  +0.37%  "

# Row 18 - Chainlink
$ws.Range("D18").Value = "12.23"
$ws.Range("E18").Value = "  +8.41%  "

# Row 19 - Polkadot
$ws.Range("E19").Value = "  +4.58%  "

# Row 20 - BitcoinCash
$ws.Range("D20").Value = "353.20"
$ws.Range("E20").Value = "  +3.88%  "

# Row 21 - Uniswap
$ws.Range("E21").Value = "  +1.65%  "

# Row 24 - Litecoin
$ws.Range("D24").Value = "66.90"
$ws.Range("E24").Value = "  +0.53%  "

# Row 25 - SuiNetwork
$ws.Range("E25").Value = "  +14.00%  "

# Row 26 - Fetch.AI
$ws.Range("D26").Value = "1.72"
$ws.Range("E26").Value = "  +6.66%  "

# Row 27 - InternetComputer(DFINITY)
$ws.Range("D27").Value = "9.41"
$ws.Range("E27").Value = "  +8.62%  "

# Row 28 - ranking swap: was Aptos, now Kaspa
$ws.Range("B28").Value = "Kaspa"
$ws.Range("C28").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D28").Value = "0.167"
$ws.Range("E28").Value = "  +2.24%  "

# Row 29 - ranking swap: was Kaspa, now Aptos
$ws.Range("B29").Value = "Aptos"
$ws.Range("C29").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D29").Value = "8.22"
$ws.Range("E29").Value = "  +3.99%  "

# Row 30 - Bittensor
$ws.Range("D30").Value = "547.57"
$ws.Range("E30").Value = "  +2.03%  "

# Row 31 - Binance-PegBSC-USD
$ws.Range("E31").Value = "  +0.19%  "

# Row 32 - PancakeSwap
$ws.Range("E32").Value = "  +1.76%  "

# Row 33 - PEPE
$ws.Range("D33").Value = "0.0₃0868"
$ws.Range("E33").Value = "  +8.04%  "

# Row 34 - ImmutableX
$ws.Range("E34").Value = "  +1.04%  "

# Row 35 - NEARProtocol
$ws.Range("E35").Value = "  +2.69%  "

# Row 36 - Monero
$ws.Range("D36").Value = "167.74"
$ws.Range("E36").Value = "  -0.55%  "

# Row 37 - Stacks
$ws.Range("D37").Value = "2.02"
$ws.Range("E37").Value = "  +7.94%  "

# Row 38 - PolygonEcosystemToken
$ws.Range("D38").Value = "0.412"
$ws.Range("E38").Value = "  +2.25%  "

# Row 39 - FirstDigitalUSD
$ws.Range("D39").Value = "0.999"
$ws.Range("E39").Value = "  -0.10%  "

# Row 40 - EthereumClassic
$ws.Range("D40").Value = "19.64"
$ws.Range("E40").Value = "  +3.37%  "

# Row 41 - USDe
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  +0.05%  "

# Row 42 - Aave
$ws.Range("D42").Value = "168.29"
$ws.Range("E42").Value = "  -0.77%  "

# Row 43 - OKB
$ws.Range("D43").Value = "40.25"
$ws.Range("E43").Value = "  +1.46%  "

# Row 44 - Filecoin
$ws.Range("D44").Value = "3.94"
$ws.Range("E44").Value = "  +5.24%  "

# Row 45 - Hedera
$ws.Range("D45").Value = "0.0587"
$ws.Range("E45").Value = "  +3.31%  "

# Row 46 - InjectiveProtocol
$ws.Range("D46").Value = "21.77"
$ws.Range("E46").Value = "  -2.81%  "

# Row 47 - ranking swap: was Mantle, now dogwifhat
$ws.Range("B47").Value = "dogwifhat"
$ws.Range("C47").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D47").Value = "2.07"
$ws.Range("E47").Value = "  +16.19%  "

# Row 48 - ranking swap: was dogwifhat, now Mantle
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").Value = "0.631"
$ws.Range("E48").Value = "  +1.19%  "

# Row 49 - VeChain
$ws.Range("E49").Value = "  +2.95%  "

# Row 50 - Stellar
$ws.Range("D50").Value = "0.0967"
$ws.Range("E50").Value = "  +0.71%  "

# Row 51 - EnergySwap
$ws.Range("D51").Value = "19.41"
$ws.Range("E51").Value = "  +5.10%  "

# Restore default styling on the Price/Volume columns so the
# temporary text NumberFormat doesn't leave a stray cell style.
$dataRange.Style = "Normal"
